$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial of 45184 for rows 2-70.
# Update it to 45185 for every such row, leaving all other values/styles intact.
$range = $ws.Range("C2:C70")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
